$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.027872562408447
$ws.Range("B1").Value = 2.58629846572876
$ws.Range("C1").Value = 2.648735761642456
$ws.Range("D1").Value = 3.541208267211914
$ws.Range("E1").Value = 5.760940074920654
